$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 19 with list name data
$ws.Range("A19").Value = "Nombre lista cliente"
$ws.Range("B19").Value = "edeq"
$ws.Range("C19").Value = "tigoUne"
$ws.Range("D19").Value = "movistar"

# Fix the IP typo in B14: "1092.168.0.1" -> "192.168.0.1"
$ws.Range("B14").Value = "192.168.0.1"

# Update the view to match the target (scroll position and selection)
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F15").Select()
